$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date for rows 2-12 from 2023-09-20 to 2023-09-21
$newDate = Get-Date -Year 2023 -Month 9 -Day 21 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
